$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "273.19"
Set-TextValue "E2" "4.49%"
Set-TextValue "D3" "26.84"
Set-TextValue "E3" "-1.13%"
Set-TextValue "D4" "4.731"
Set-TextValue "E4" "0.59%"
Set-TextValue "D5" "0.06203"
Set-TextValue "E5" "-0.13%"
Set-TextValue "D6" "6.769"
Set-TextValue "E6" "0.65%"
Set-TextValue "D7" "0.8628"
Set-TextValue "E7" "1.37%"
Set-TextValue "D8" "0.9103"
Set-TextValue "E8" "0.41%"
Set-TextValue "D9" "0.1439"
Set-TextValue "E9" "2.64%"
Set-TextValue "D10" "0.05360"
Set-TextValue "E10" "13.92%"
Set-TextValue "D11" "0.07184"
Set-TextValue "E11" "1.35%"
Set-TextValue "D12" "0.03182"
Set-TextValue "E12" "0.28%"
Set-TextValue "D13" "0.09051"
Set-TextValue "E13" "-0.02%"
Set-TextValue "D14" "0.001534"
Set-TextValue "E14" "-0.74%"
Set-TextValue "D15" "0.0006084"
Set-TextValue "E15" "-0.93%"
Set-TextValue "D16" "0.005937"
Set-TextValue "E16" "-1.09%"
Set-TextValue "D17" "3.470"
Set-TextValue "E17" "0.10%"
Set-TextValue "D18" "3.194"
Set-TextValue "E18" "0.73%"
Set-TextValue "E19" "4.04%"
Set-TextValue "E21" "1.36%"
Set-TextValue "D22" "3.845"
Set-TextValue "E22" "-6.18%"
Set-TextValue "E23" "0.51%"
Set-TextValue "D24" "0.001175"
Set-TextValue "E24" "-3.67%"
Set-TextValue "D25" "0.004195"
Set-TextValue "E25" "1.98%"
Set-TextValue "D26" "0.0001198"
Set-TextValue "E26" "-0.20%"
Set-TextValue "E40" "1.92%"
Set-TextValue "D41" "0.006216"
Set-TextValue "E41" "50.34%"
Set-TextValue "E42" "1.54%"
Set-TextValue "D43" "0.002169"
Set-TextValue "E43" "-0.64%"
Set-TextValue "D44" "0.01282"
Set-TextValue "E44" "-4.46%"
Set-TextValue "D45" "0.00005130"
Set-TextValue "E45" "-0.22%"
Set-TextValue "E46" "-0.24%"
Set-TextValue "D47" "0.8954"
Set-TextValue "E47" "3,282.59%"
Set-TextValue "E48" "-14.68%"
Set-TextValue "D49" "0.00002096"
Set-TextValue "E49" "-0.24%"
Set-TextValue "D50" "0.0001997"
Set-TextValue "E50" "-0.24%"
